$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Modelo" in F1, reusing the same formatting as the other headers (A1:E1)
$ws.Range("F1").Value = "Modelo"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Update numeric values in row 2 (MSE, R2, MAE)
$ws.Range("B2").Value = 0.0675542386722679
$ws.Range("C2").Value = 0.9987659440417859
$ws.Range("D2").Value = 0.2056142416656275

# Add model description text in F2
$ws.Range("F2").Value = "Pipeline(steps=[('model',`n                 RandomForestRegressor(max_depth=5, n_estimators=150))])"
